# Adapt column header formatting to respective input-file names:
#   *_old -> *_FV2410   (columns A..J)
#   *_new -> *_FV2504   (columns L..U)
# then wrap the used range in a table and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Rename the header cells (row 1).
# ---------------------------------------------------------------------------
$oldSuffixCols = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$newSuffixCols = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

foreach ($col in $oldSuffixCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Value2 -replace "_old$", "_FV2410")
}

foreach ($col in $newSuffixCols) {
    $cell = $ws.Range($col + "1")
    $cell.Value = ($cell.Value2 -replace "_new$", "_FV2504")
}

# ---------------------------------------------------------------------------
# 2. Turn the used range into a native Excel table ("Table1") without
#    picking up an automatic banded/default style dxf on top of the
#    pre-existing header formatting (bold + fill + border already baked
#    into the sheet). We stash the header formatting away, clear it so
#    table-creation has nothing special to "absorb" into a header dxf,
#    then restore it once the table exists.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:U1")
$scratchRange = $ws.Range("A1000:U1000")

$headerRange.Copy()
$scratchRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$headerRange.ClearFormats()

$tableRange = $ws.Range("A1:U70")
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $tableRange, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"
$lo.TableStyle = ""

$scratchRange.Copy()
$headerRange.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$scratchRange.ClearFormats()
$scratchRange.ClearContents()

$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# 3. Freeze the header row (pane split after row 1).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
